$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 0.2917716402565462;  C = 0.306821227259698;   D = 0.1494219747398047; E = 0.4942365360607697;  G = 1.242251378316819 }
    3 = @{ B = 0.6606524410359556;  C = 1.655778082260271;   D = 0.1494219747398047; E = 0.4942365360607697;  G = 2.960089034096801 }
    4 = @{ B = 0.1190320826869504;  C = 0.002571899574220771; D = 0.1494219747398047; E = 0.4942365360607697;  G = 0.7652624930617455 }
    5 = @{ B = 0.6606524410359556;  C = 1.655778082260271;   D = 0.1494219747398047; E = 10.19245300693656;   G = 12.65830550497259 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
